$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date for every data row
# (rows 2-530). The update bumps this date by one day (2023-09-02 ->
# 2023-09-03, serial 45171 -> 45172) for every row.
$ws.Range("C2:C530").Value = 45172
